# Updates cryptos list prices/volumes (and the Aptos/ICP row swap) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.138.70'
$ws.Range('E2').Value = '  +0.34%  '

# Row 3
$ws.Range('D3').Value = '2.604.88'
$ws.Range('E3').Value = '  +2.26%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').Value = '''583.48'
$ws.Range('E5').Value = '  +2.50%  '

# Row 6
$ws.Range('D6').Value = '''148.06'
$ws.Range('E6').Value = '  +1.26%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').Value = '''0.599'
$ws.Range('E8').Value = '  +2.75%  '

# Row 9
$ws.Range('E9').Value = '  +3.20%  '

# Row 10
$ws.Range('D10').Value = '''5.69'
$ws.Range('E10').Value = '  +3.29%  '

# Row 11
$ws.Range('E11').Value = '  +0.08%  '

# Row 12
$ws.Range('D12').Value = '''0.355'
$ws.Range('E12').Value = '  +0.51%  '

# Row 13
$ws.Range('D13').Value = '''27.31'
$ws.Range('E13').Value = '  -0.08%  '

# Row 14
$ws.Range('D14').Value = '3.072.70'
$ws.Range('E14').Value = '  +2.34%  '

# Row 15
$ws.Range('D15').Value = '63.051.09'
$ws.Range('E15').Value = '  +0.25%  '

# Row 16
$ws.Range('D16').Value = '''0.0000148'
$ws.Range('E16').Value = '  +3.23%  '

# Row 17
$ws.Range('D17').Value = '2.603.09'

# Row 18
$ws.Range('D18').Value = '''11.38'
$ws.Range('E18').Value = '  +0.69%  '

# Row 19
$ws.Range('D19').Value = '''343.96'
$ws.Range('E19').Value = '  +2.76%  '

# Row 20
$ws.Range('D20').Value = '''4.42'
$ws.Range('E20').Value = '  +1.64%  '

# Row 21
$ws.Range('D21').Value = '''6.79'
$ws.Range('E21').Value = '  +0.12%  '

# Row 22
$ws.Range('E22').Value = '  +0.03%  '

# Row 23
$ws.Range('D23').Value = '''5.66'
$ws.Range('E23').Value = '  -1.60%  '

# Row 24
$ws.Range('D24').Value = '''66.98'
$ws.Range('E24').Value = '  +2.61%  '

# Row 25
$ws.Range('D25').Value = '2.723.25'
$ws.Range('E25').Value = '  +1.72%  '

# Row 26
$ws.Range('E26').Value = '  +0.08%  '

# Row 27
$ws.Range('E27').Value = '  +0.83%  '

# Row 28
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  +0.02%  '

# Row 29
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''8.44'
$ws.Range('E29').Value = '  +1.09%  '

# Row 30
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '''7.90'
$ws.Range('E30').Value = '  +7.84%  '

# Row 31
$ws.Range('D31').Value = '''1.46'
$ws.Range('E31').Value = '  -0.43%  '

# Row 32
$ws.Range('D32').Value = '''1.94'
$ws.Range('E32').Value = '  +5.25%  '

# Row 33
$ws.Range('D33').Value = '0.0₃0824'
$ws.Range('E33').Value = '  +1.35%  '

# Row 34
$ws.Range('D34').Value = '''465.47'
$ws.Range('E34').Value = '  +14.79%  '

# Row 35
$ws.Range('D35').Value = '''176.87'
$ws.Range('E35').Value = '  +0.71%  '

# Row 36
$ws.Range('E36').Value = '  +4.46%  '

# Row 37
$ws.Range('E37').Value = '  +0.06%  '

# Row 38
$ws.Range('D38').Value = '''0.403'
$ws.Range('E38').Value = '  +0.60%  '

# Row 39
$ws.Range('D39').Value = '''19.26'
$ws.Range('E39').Value = '  +0.88%  '

# Row 40
$ws.Range('D40').Value = '''4.61'
$ws.Range('E40').Value = '  +6.27%  '

# Row 41
$ws.Range('E41').Value = '  +0.03%  '

# Row 42
$ws.Range('D42').Value = '''1.71'
$ws.Range('E42').Value = '  -1.82%  '

# Row 43
$ws.Range('D43').Value = '''160.51'
$ws.Range('E43').Value = '  +5.63%  '

# Row 44
$ws.Range('D44').Value = '''3.80'
$ws.Range('E44').Value = '  +1.25%  '

# Row 45
$ws.Range('D45').Value = '''0.639'
$ws.Range('E45').Value = '  +6.28%  '

# Row 46
$ws.Range('D46').Value = '''20.92'
$ws.Range('E46').Value = '  +0.45%  '

# Row 47
$ws.Range('D47').Value = '''0.0548'
$ws.Range('E47').Value = '  +2.87%  '

# Row 48
$ws.Range('D48').Value = '''0.0976'
$ws.Range('E48').Value = '  +1.13%  '

# Row 49
$ws.Range('D49').Value = '''0.0238'
$ws.Range('E49').Value = '  -0.42%  '

# Row 50
$ws.Range('D50').Value = '''18.70'
$ws.Range('E50').Value = '  +2.25%  '

# Row 51
$ws.Range('E51').Value = '  -0.03%  '
